# Auto-generated script applying scheduled-runner price/profit updates
# to the Asura_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2173.25
$ws.Range("I15").Value = 2173.25
$ws.Range("K15").Value = 6519.75
$ws.Range("M15").Value = -6350.75
$ws.Range("H17").Value = 460.8
$ws.Range("J17").Value = 460.8
$ws.Range("L17").Value = 1382.4
$ws.Range("N17").Value = -1718.4
$ws.Range("H43").Value = 1078
$ws.Range("I43").Value = 1119.2222
$ws.Range("J43").Value = 985.25
$ws.Range("K43").Value = 1119.2222
$ws.Range("L43").Value = 985.25
$ws.Range("M43").Value = -1050.2222
$ws.Range("N43").Value = -1123.25
$ws.Range("H55").Value = 562.5
$ws.Range("I55").Value = 200.2
$ws.Range("J55").Value = 1166.3334
$ws.Range("K55").Value = 200.2
$ws.Range("L55").Value = 1166.3334
$ws.Range("M55").Value = 13.80000000000001
$ws.Range("N55").Value = -1594.3334
$ws.Range("H88").Value = 3019.889
$ws.Range("I88").Value = 2779
$ws.Range("J88").Value = 3050
$ws.Range("K88").Value = 2779
$ws.Range("L88").Value = 3050
$ws.Range("M88").Value = -2373
$ws.Range("N88").Value = -3862
$ws.Range("H91").Value = 3019.889
$ws.Range("I91").Value = 2779
$ws.Range("J91").Value = 3050
$ws.Range("K91").Value = 2779
$ws.Range("L91").Value = 3050
$ws.Range("M91").Value = -1375
$ws.Range("N91").Value = -5858
$ws.Range("H106").Value = 3170.111
$ws.Range("I106").Value = 3937.5
$ws.Range("K106").Value = 3937.5
$ws.Range("M106").Value = -3306.5
$ws.Range("H116").Value = 9093322
$ws.Range("I116").Value = 66667436
$ws.Range("J116").Value = 2672.6316
$ws.Range("K116").Value = 66667436
$ws.Range("L116").Value = 2672.6316
$ws.Range("M116").Value = -66663994
$ws.Range("N116").Value = -9556.631600000001
$ws.Range("H129").Value = 1139.2826
$ws.Range("I129").Value = 472.44446
$ws.Range("J129").Value = 1301.4865
$ws.Range("K129").Value = 1417.33338
$ws.Range("L129").Value = 3904.4595
$ws.Range("M129").Value = 3582.66662
$ws.Range("N129").Value = -13904.4595
$ws.Range("H137").Value = 1327.3513
$ws.Range("I137").Value = 1314.8276
$ws.Range("J137").Value = 1372.75
$ws.Range("K137").Value = 3944.4828
$ws.Range("L137").Value = 4118.25
$ws.Range("M137").Value = -1394.4828
$ws.Range("N137").Value = -9218.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9635.529
$ws.Range("I32").Value = 10461.637
$ws.Range("J32").Value = 6140.4614
$ws.Range("K32").Value = 10461.637
$ws.Range("L32").Value = 6140.4614
$ws.Range("M32").Value = -10174.637
$ws.Range("N32").Value = -6714.4614
$ws.Range("H61").Value = 2136.4285
$ws.Range("I61").Value = 1900.909
$ws.Range("K61").Value = 1900.909
$ws.Range("M61").Value = -1688.909
$ws.Range("H74").Value = 961.48834
$ws.Range("I74").Value = 833.91174
$ws.Range("K74").Value = 833.91174
$ws.Range("M74").Value = 40.08825999999999
$ws.Range("H77").Value = 961.48834
$ws.Range("I77").Value = 833.91174
$ws.Range("K77").Value = 4169.5587
$ws.Range("M77").Value = 198.4413000000004
$ws.Range("H132").Value = 4461.0654
$ws.Range("I132").Value = 5141.9033
$ws.Range("J132").Value = 3054
$ws.Range("K132").Value = 15425.7099
$ws.Range("L132").Value = 9162
$ws.Range("M132").Value = -12895.7099
$ws.Range("N132").Value = -14222
$ws.Range("H136").Value = 2136.4285
$ws.Range("I136").Value = 1900.909
$ws.Range("K136").Value = 5702.727000000001
$ws.Range("M136").Value = -3152.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 6713.857
$ws.Range("I5").Value = 6713.857
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6713.857
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -6600.857
$ws.Range("N5").ClearContents()
$ws.Range("H132").Value = 84253.336
$ws.Range("J132").Value = 84253.336
$ws.Range("L132").Value = 84253.336
$ws.Range("N132").Value = -94373.336
$ws.Range("H134").Value = 2193.827
$ws.Range("I134").Value = 1874.1395
$ws.Range("J134").Value = 3721.2222
$ws.Range("K134").Value = 5622.4185
$ws.Range("L134").Value = 11163.6666
$ws.Range("M134").Value = -3087.4185
$ws.Range("N134").Value = -16233.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8857.143
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H31").Value = 2293.7856
$ws.Range("I31").Value = 1296.1052
$ws.Range("K31").Value = 1296.1052
$ws.Range("M31").Value = -1001.1052
$ws.Range("H34").Value = 2293.7856
$ws.Range("I34").Value = 1296.1052
$ws.Range("K34").Value = 1296.1052
$ws.Range("M34").Value = -1094.1052
$ws.Range("H58").Value = 863199.4399999999
$ws.Range("I58").Value = 1196385.9
$ws.Range("J58").Value = 2467.8333
$ws.Range("K58").Value = 1196385.9
$ws.Range("L58").Value = 2467.8333
$ws.Range("M58").Value = -1196182.9
$ws.Range("N58").Value = -2873.8333
$ws.Range("H122").Value = 1680.9459
$ws.Range("I122").Value = 1708.3043
$ws.Range("J122").Value = 1636
$ws.Range("K122").Value = 5124.9129
$ws.Range("L122").Value = 4908
$ws.Range("M122").Value = -2674.9129
$ws.Range("N122").Value = -9808
$ws.Range("H132").Value = 521649.72
$ws.Range("I132").Value = 615821.9
$ws.Range("J132").Value = 3703
$ws.Range("K132").Value = 1847465.7
$ws.Range("L132").Value = 11109
$ws.Range("M132").Value = -1844935.7
$ws.Range("N132").Value = -16169
$ws.Range("H134").Value = 1566.8928
$ws.Range("I134").Value = 972.13043
$ws.Range("J134").Value = 4302.8
$ws.Range("K134").Value = 2916.39129
$ws.Range("L134").Value = 12908.4
$ws.Range("M134").Value = -381.39129
$ws.Range("N134").Value = -17978.4
$ws.Range("H136").Value = 863199.4399999999
$ws.Range("I136").Value = 1196385.9
$ws.Range("J136").Value = 2467.8333
$ws.Range("K136").Value = 3589157.7
$ws.Range("L136").Value = 7403.499899999999
$ws.Range("M136").Value = -3586607.7
$ws.Range("N136").Value = -12503.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 614.6667
$ws.Range("I113").Value = 529.41174
$ws.Range("K113").Value = 1588.23522
$ws.Range("M113").Value = 581.76478
$ws.Range("H132").Value = 1763.1428
$ws.Range("I132").Value = 1263.1818
$ws.Range("K132").Value = 11368.6362
$ws.Range("M132").Value = -8838.636200000001
$ws.Range("H133").Value = 3802.087
$ws.Range("I133").Value = 1813.0769
$ws.Range("K133").Value = 5439.2307
$ws.Range("M133").Value = -379.2307000000001
$ws.Range("H136").Value = 3866.6667
$ws.Range("I136").Value = 921.1818
$ws.Range("K136").Value = 2763.5454
$ws.Range("M136").Value = 2336.4546
$ws.Range("H137").Value = 9011301
$ws.Range("I137").Value = 998.125
$ws.Range("J137").Value = 15876294
$ws.Range("K137").Value = 2994.375
$ws.Range("L137").Value = 47628882
$ws.Range("M137").Value = 2105.625
$ws.Range("N137").Value = -47639082

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2359.3928
$ws.Range("I102").Value = 1797.65
$ws.Range("J102").Value = 3763.75
$ws.Range("K102").Value = 1797.65
$ws.Range("L102").Value = 3763.75
$ws.Range("M102").Value = -175.6500000000001
$ws.Range("N102").Value = -7007.75
$ws.Range("H132").Value = 2629.1
$ws.Range("I132").Value = 1874.85
$ws.Range("J132").Value = 4137.6
$ws.Range("K132").Value = 5624.549999999999
$ws.Range("L132").Value = 12412.8
$ws.Range("M132").Value = -3094.549999999999
$ws.Range("N132").Value = -17472.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1025.0625
$ws.Range("I46").Value = 781.9091
$ws.Range("K46").Value = 781.9091
$ws.Range("M46").Value = -593.9091
$ws.Range("H61").Value = 33990.332
$ws.Range("I61").Value = 40596.4
$ws.Range("J61").Value = 960
$ws.Range("K61").Value = 40596.4
$ws.Range("L61").Value = 960
$ws.Range("M61").Value = -40394.4
$ws.Range("N61").Value = -1364
$ws.Range("H68").Value = 2700
$ws.Range("I68").Value = 2530.7693
$ws.Range("J68").Value = 3250
$ws.Range("K68").Value = 2530.7693
$ws.Range("L68").Value = 3250
$ws.Range("M68").Value = -1781.7693
$ws.Range("N68").Value = -4748
$ws.Range("H71").Value = 2700
$ws.Range("I71").Value = 2530.7693
$ws.Range("J71").Value = 3250
$ws.Range("K71").Value = 12653.8465
$ws.Range("L71").Value = 16250
$ws.Range("M71").Value = -8909.8465
$ws.Range("N71").Value = -23738
$ws.Range("H113").Value = 33990.332
$ws.Range("I113").Value = 40596.4
$ws.Range("J113").Value = 960
$ws.Range("K113").Value = 40596.4
$ws.Range("L113").Value = 960
$ws.Range("M113").Value = -38426.4
$ws.Range("N113").Value = -5300
$ws.Range("H132").Value = 3914.1482
$ws.Range("I132").Value = 3900.2856
$ws.Range("K132").Value = 11700.8568
$ws.Range("M132").Value = -9170.856800000001
$ws.Range("H136").Value = 27299692
$ws.Range("I136").Value = 47620644
$ws.Range("J136").Value = 628444.4
$ws.Range("K136").Value = 142861932
$ws.Range("L136").Value = 1885333.2
$ws.Range("M136").Value = -142859382
$ws.Range("N136").Value = -1890433.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1685.902
$ws.Range("I132").Value = 1105.7188
$ws.Range("J132").Value = 2663.0527
$ws.Range("K132").Value = 3317.1564
$ws.Range("L132").Value = 7989.158100000001
$ws.Range("M132").Value = -787.1564000000003
$ws.Range("N132").Value = -13049.1581
$ws.Range("H136").Value = 1428.8788
$ws.Range("I136").Value = 1190.8518
$ws.Range("K136").Value = 3572.5554
$ws.Range("M136").Value = -1022.5554
